$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Design and implementation of a Contact Center ... with a RESTful API, "
#    -> insert "Suite " right before "with a RESTful API, "
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("with a RESTful API, ", $true, $false, $false, $false, $false, $true, 1, $false, `
    "Suite with a RESTful API, ", 2)
Write-Output "Suite insert found=$found"

# ---------------------------------------------------------------------------
# 2) OpenALPR bullet: collapse the three runs ("...version ", "working with",
#    the _GoBack bookmark, " Venezuela") that spelled out
#    "OpenALPR (Automatic License Plate Recognition) version working with Venezuela"
#    into one run and drop the stray _GoBack bookmark that used to sit in the
#    middle of it (it gets re-created further down, next to the new R2C2 text).
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute( `
    "OpenALPR (Automatic License Plate Recognition) version working with Venezuela", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "OpenALPR (Automatic License Plate Recognition) version working with Venezuela", 2)
Write-Output "OpenALPR normalize found=$found2"

# ---------------------------------------------------------------------------
# 3) City route bullet -> add the R2C2 4500x-speedup blurb.
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$found3 = $rng3.Find.Execute( `
    "City route recommendation system based on Tweets about vehicular transit. Working in real time.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "4500x improvement to a city route recommendation system based on Tweets about vehicular transit. Working in real time thanks to this optimization.", 2)
Write-Output "City route rewrite found=$found3"

# Re-plant the _GoBack bookmark right before "recommendation system...",
# matching where Word would leave the caret after that edit.
$rng4 = $d.Content
$found4 = $rng4.Find.Execute("recommendation system based on Tweets", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Output "bookmark anchor found=$found4 start=$($rng4.Start)"
if ($found4) {
    $bmRange = $d.Range($rng4.Start, $rng4.Start)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

# ---------------------------------------------------------------------------
# 4) Row-height relayout for the Coursera / Stanford rows caused by the new
#    text above reflowing the table.
# ---------------------------------------------------------------------------
$t = $d.Tables.Item(1)
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $row = $t.Rows.Item($i)
    $cellText = $row.Cells.Item(1).Range.Text
    if ($cellText -like "Coursera*") {
        $row.Height = 22.4
        Write-Output "Set Coursera row ($i) height -> $($row.Height)"
    } elseif ($cellText -like "Stanford University*") {
        $row.Height = 95.9
        Write-Output "Set Stanford row ($i) height -> $($row.Height)"
    }
}
